# Add a new Job Posting row (Job_Id = 11) to the LinkedIn job posting sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 12

$ws.Cells.Item($newRow, 1).Value = 11                # Job_Id
$ws.Cells.Item($newRow, 2).Value = "React Developer"  # Jd_Title (reuses existing shared string)
$ws.Cells.Item($newRow, 3).Value = "fafwasf"          # Job_Description (new shared string)
$ws.Cells.Item($newRow, 4).Value = 1                  # Total_Years_Min_Exp
$ws.Cells.Item($newRow, 5).Value = 3                  # Total_Years_Max_Exp
$ws.Cells.Item($newRow, 6).Value = 0                  # Linked_Posted
$ws.Cells.Item($newRow, 7).Value = 0                  # Resume_received
$ws.Cells.Item($newRow, 8).Value = 0                  # Resume_downloaded
